# Apply the "Updated test data for normal load, cable capacitance etc" edit
# to the "Add Devices" sheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "Add Devices"

# --- Update the "Other Slot Cards (... of ...)" labels (typo: missing closing paren) ---
$ws.Range("O10").Value = "Other Slot Cards  (2 of 3"
$ws.Range("O11").Value = "Other Slot Cards  (2 of 4"
$ws.Range("O12").Value = "Other Slot Cards  (1 of 4"
$ws.Range("O13").Value = "Other Slot Cards  (1 of 4"

# --- Update column headers ---
$ws.Range("K8").Value = "IOB800(x2)"
$ws.Range("K9").Value = "AttachedFunctionality"

# --- Re-set B3 so it reuses the existing shared string (index shifts after cleanup) ---
$ws.Range("B3").Value = "VerifyAddUnitDetails"

# --- Update numeric data values ---
$ws.Range("K11").Value = 1
$ws.Range("H13").Value = 1
$ws.Range("I13").Value = 0

# --- Update the active selection on the sheet ---
$ws.Range("H11").Select() | Out-Null
